$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Investment Type *" column (old column F). Rather than deleting
# the whole column outright, shift the columns to its right (Notes, Type,
# Folio No) one place to the left and clear the now-unused column I.
for ($r = 1; $r -le 9; $r++) {
    $g = $ws.Cells.Item($r, 7).Value()
    $h = $ws.Cells.Item($r, 8).Value()
    $i = $ws.Cells.Item($r, 9).Value()

    $ws.Cells.Item($r, 6).Value = $g
    $ws.Cells.Item($r, 7).Value = $h
    $ws.Cells.Item($r, 8).Value = $i
    $ws.Cells.Item($r, 9).Value = $null
}

# New header: Sector / Category / Sub Category / Startup / Investment Domicile
$ws.Range("I1").Value = "Sector"
$ws.Range("J1").Value = "Category"
$ws.Range("K1").Value = "Sub Category"
$ws.Range("L1").Value = "Startup"
$ws.Range("M1").Value = "Investment Domicile"

# Fill in the new Category / Sub Category / Startup / Investment Domicile
# values for every data row (2-9). Sector (column I) stays blank.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 10).Value = "Unlisted"
    $ws.Cells.Item($r, 11).Value = "Equity"
    $ws.Cells.Item($r, 12).Value = "Yes"
    $ws.Cells.Item($r, 13).Value = "Domestic"
}

$ws.Range("J3:J9").Select()
